{"js": "// Word JS API (Office.js) edit script\n// Body of: async (context) => { ... }\n//\n// Changes applied (per the supplied OOXML diff):\n//   1. Queue section: both \"<quote> ben\u00f6tigt.\" sentences gain a trailing\n//      qualifier, becoming \"<quote> ben\u00f6tigt (standardm\u00e4ssig in\n//      standalone-full enthalten).\" \u2014 for the \"...ConnectionFactory\"\n//      paragraph and the \"...queue/test\" paragraph.\n//   2. The Word-managed \"_GoBack\" bookmark (marks the last edit location)\n//      moves from the \"sql scripts\" folder description paragraph to the\n//      start of the \"Verwendete Technologien\" heading paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst suffix = \" (standardm\u00e4ssig in standalone-full enthalten)\";\n\n// Locate, inside each paragraph, the literal run of text that ends the\n// sentence (\"\u2026 ben\u00f6tigt.\") and splice in the new qualifier just before the\n// final period, leaving everything else (including closing quote) intact.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text || \"\";\n  if (text.indexOf(\"ConnectionFactory\") !== -1 && text.trim().endsWith(\"ben\u00f6tigt.\")) {\n    const hits = para.search(\"\u201c ben\u00f6tigt.\", { matchCase: true });\n    hits.load(\"items/text\");\n    await context.sync();\n    if (hits.items.length > 0) {\n      hits.items[0].insertText(\"\u201c ben\u00f6tigt\" + suffix + \".\", \"Replace\");\n      await context.sync();\n    }\n  } else if (text.indexOf(\"queue/test\") !== -1 && text.trim().endsWith(\"ben\u00f6tigt.\")) {\n    const hits = para.search(\"\u201c ben\u00f6tigt.\", { matchCase: true });\n    hits.load(\"items/text\");\n    await context.sync();\n    if (hits.items.length > 0) {\n      hits.items[0].insertText(\"\u201c ben\u00f6tigt\" + suffix + \".\", \"Replace\");\n      await context.sync();\n    }\n  }\n}\n\n// Move the \"_GoBack\" bookmark to the start of the \"Verwendete Technologien\"\n// heading paragraph (removing it from wherever it currently sits).\nconst existing = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexisting.load(\"isNullObject\");\nawait context.sync();\nif (!existing.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nconst headingHits = context.document.body.search(\"Verwendete Technologien\", { matchCase: true });\nheadingHits.load(\"items\");\nawait context.sync();\nif (headingHits.items.length > 0) {\n  const headingRange = headingHits.items[0].getRange(\"Start\");\n  headingRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script\n# $word.ActiveDocument is the open document.\n#\n# Changes applied (per the supplied OOXML diff):\n#   1. Queue section: both \"<quote> ben\u00f6tigt.\" sentences gain a trailing\n#      qualifier, becoming \"<quote> ben\u00f6tigt (standardm\u00e4ssig in\n#      standalone-full enthalten).\" \u2014 for the \"...ConnectionFactory\"\n#      paragraph and the \"...queue/test\" paragraph.\n#   2. The Word-managed \"_GoBack\" bookmark (marks the last edit location)\n#      moves from the \"sql scripts\" folder description paragraph to the\n#      start of the \"Verwendete Technologien\" heading paragraph.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n$replacement = \"ben\u00f6tigt (standardm\u00e4ssig in standalone-full enthalten).\"\n\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*ConnectionFactory*ben\u00f6tigt.*\" -or $t -like \"*queue/test*ben\u00f6tigt.*\") {\n        $r = $p.Range.Duplicate\n        $r.Find.ClearFormatting()\n        $r.Find.Text = \"ben\u00f6tigt.\"\n        $r.Find.Replacement.ClearFormatting()\n        $r.Find.Replacement.Text = $replacement\n        $r.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n    }\n}\n\n# Move the \"_GoBack\" bookmark to the start of the \"Verwendete Technologien\"\n# heading paragraph (removing it from wherever it currently sits).\ntry {\n    $bm = $d.Bookmarks.Item(\"_GoBack\")\n    $bm.Delete()\n} catch {\n}\n\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"Verwendete Technologien*\") {\n        $r = $p.Range.Duplicate\n        $r.Collapse(1)\n        $d.Bookmarks.Add(\"_GoBack\", $r)\n        break\n    }\n}\n"}
